$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['SoftwareFault']"

$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

$ws.Range("D56").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['Normal']"
